$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.552.54"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "1.849.75"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("D4").Value = "'0.9992"

$ws.Range("D5").Value = "'240.56"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "'0.6298"
$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "'0.07463"
$ws.Range("E8").Value = "  -1.66%  "

$ws.Range("D9").Value = "'0.2911"
$ws.Range("E9").Value = "  -0.34%  "

$ws.Range("D10").Value = "'25.01"
$ws.Range("E10").Value = "  +1.96%  "

$ws.Range("D11").Value = "'0.07745"
$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("D12").Value = "1.855.12"
$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("D13").Value = "'5.018"
$ws.Range("E13").Value = "  +0.17%  "

$ws.Range("D14").Value = "'0.6832"
$ws.Range("E14").Value = "  +0.63%  "

$ws.Range("D15").Value = "'0.00001024"
$ws.Range("E15").Value = "  -1.76%  "

$ws.Range("D16").Value = "'82.62"
$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("D17").Value = "'6.322"
$ws.Range("E17").Value = "  +3.51%  "

$ws.Range("D18").Value = "29.566.16"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").Value = "'230.07"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").Value = "'12.38"
$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("D21").Value = "'0.9996"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").Value = "'7.513"
$ws.Range("E22").Value = "  +1.11%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'159.22"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").Value = "'8.510"
$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("D26").Value = "'0.1365"
$ws.Range("E26").Value = "  -2.17%  "

$ws.Range("D27").Value = "'17.57"
$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("B28").Value = "Hedera"
$ws.Range("C28").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D28").Value = "'0.06611"
$ws.Range("E28").Value = "  +16.27%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.470"
$ws.Range("E29").Value = "  +3.15%  "

$ws.Range("D30").Value = "'1.487"
$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("D31").Value = "'4.105"
$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("D32").Value = "'4.091"
$ws.Range("E32").Value = "  +1.19%  "

$ws.Range("D33").Value = "'1.851"
$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("D34").Value = "'1.141"
$ws.Range("E34").Value = "  -1.29%  "

$ws.Range("D35").Value = "'0.6985"
$ws.Range("E35").Value = "  +0.36%  "

$ws.Range("D36").Value = "'2.563"
$ws.Range("E36").Value = "  -0.78%  "

$ws.Range("E37").Value = "  +2.37%  "

$ws.Range("D38").Value = "'2.839"
$ws.Range("E38").Value = "  +4.43%  "

$ws.Range("D39").Value = "1.254.10"
$ws.Range("E39").Value = "  +1.19%  "

$ws.Range("D40").Value = "'6.782"
$ws.Range("E40").Value = "  +5.73%  "

$ws.Range("D41").Value = "'0.9352"
$ws.Range("E41").Value = "  +3.69%  "

$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("D43").Value = "2.008.60"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").Value = "'101.42"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").Value = "'66.26"
$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("D46").Value = "'7.099"
$ws.Range("E46").Value = "  -0.53%  "

$ws.Range("D47").Value = "'1.729"
$ws.Range("E47").Value = "  +2.95%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.061"
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.1157"
$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000115"
$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "'0.00000000115"
$ws.Range("E51").Value = "  -1.34%  "

